$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel keeps them as strings
# (matches the source data which stores prices as text, e.g. thousands-dot formatted).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "62.120.69"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.424.97"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "563.95"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "144.09"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "2.423.98"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "26.22"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "2.859.28"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "61.979.26"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.433.73"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "324.23"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "67.31"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").Value = "559.55"
$ws.Range("E27").Value = "  -5.13%  "
$ws.Range("D28").Value = "2.541.45"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").Value = "8.23"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -5.34%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "5.49"
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("D40").Value = "152.31"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").Value = "18.69"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "147.85"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "0.0532"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").Value = "19.96"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "0.597"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -0.52%  "
